$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    4  = @{ C = -0.5865396705125869;  E = 0.8461488352962343 }
    5  = @{ C = 2.46453841416181;     E = 1.283106670049072 }
    6  = @{ C = 1.024488446359717;    E = 0.883996728309433 }
    7  = @{ C = -0.2650887206215868;  E = 0.7934423576064997 }
    8  = @{ C = 2.011885157458759;    E = 1.708833114833652 }
    9  = @{ C = 1.582403480979067;    E = 1.315455703544299 }
    10 = @{ C = 2.232192661479382;    E = 1.471326103602411 }
    11 = @{ C = 1.672539262716755;    E = 1.230669741902624 }
    12 = @{ C = 1.547454231976442;    E = 1.200013281595735 }
    13 = @{ C = 1.380811415368788;    E = 1.063817241834708 }
    14 = @{ C = -0.9147268599428826;  E = 0.1022244835662045 }
    15 = @{ C = -0.6863413708432242;  E = 1.199013851879926 }
    16 = @{ C = 3.676169537136298;    E = 1.620572283829613 }
    17 = @{ C = -0.3214161021417694;  E = 0.7147372843446353 }
    18 = @{ C = -0.968109702189679;   E = 0.7595446561616592 }
    19 = @{ C = 1.686416487556031;    E = 0.8971523338365817 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
